# Auto-generated edit script: update FFXIV leve-profit market-data sheets
# per scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")

# --- ALC ---
# row 51
$ws_ALC.Range("H51").Value = 5540
$ws_ALC.Range("I51").Value = 2350
$ws_ALC.Range("J51").Value = 7666.6665
$ws_ALC.Range("K51").Value = 2350
$ws_ALC.Range("L51").Value = 7666.6665
$ws_ALC.Range("M51").Value = -1866
$ws_ALC.Range("N51").Value = -8634.666499999999

# row 55
$ws_ALC.Range("H55").Value = 270.86667
$ws_ALC.Range("I55").Value = 357
$ws_ALC.Range("J55").Value = 141.66667
$ws_ALC.Range("K55").Value = 357
$ws_ALC.Range("L55").Value = 141.66667
$ws_ALC.Range("M55").Value = -143
$ws_ALC.Range("N55").Value = -569.6666700000001

# row 100
$ws_ALC.Range("H100").Value = 6667930.5
$ws_ALC.Range("I100").Value = 9805022
$ws_ALC.Range("J100").Value = 1612
$ws_ALC.Range("K100").Value = 9805022
$ws_ALC.Range("L100").Value = 1612
$ws_ALC.Range("M100").Value = -9804481
$ws_ALC.Range("N100").Value = -2694

# row 113
$ws_ALC.Range("H113").Value = 4953.1113
$ws_ALC.Range("I113").Value = 2915.6
$ws_ALC.Range("J113").Value = 7500
$ws_ALC.Range("K113").Value = 2915.6
$ws_ALC.Range("L113").Value = 7500
$ws_ALC.Range("M113").Value = 338.4000000000001
$ws_ALC.Range("N113").Value = -14008

# row 132
$ws_ALC.Range("H132").Value = 2165.2334
$ws_ALC.Range("I132").Value = 1275.1154
$ws_ALC.Range("K132").Value = 3825.3462
$ws_ALC.Range("M132").Value = -1295.3462

# --- ARM ---
# row 6
$ws_ARM.Range("H6").Value = 410251.5
$ws_ARM.Range("I6").Value = 1000
$ws_ARM.Range("J6").Value = 546668.7
$ws_ARM.Range("K6").Value = 1000
$ws_ARM.Range("L6").Value = 546668.7
$ws_ARM.Range("M6").Value = -827
$ws_ARM.Range("N6").Value = -547014.7

# row 32
$ws_ARM.Range("H32").Value = 4729.2817
$ws_ARM.Range("I32").Value = 3544.0356
$ws_ARM.Range("J32").Value = 9154.200000000001
$ws_ARM.Range("K32").Value = 3544.0356
$ws_ARM.Range("L32").Value = 9154.200000000001
$ws_ARM.Range("M32").Value = -3257.0356
$ws_ARM.Range("N32").Value = -9728.200000000001

# row 45
$ws_ARM.Range("H45").Value = 6916.2104
$ws_ARM.Range("I45").Value = 25931.25
$ws_ARM.Range("J45").Value = 1845.5333
$ws_ARM.Range("K45").Value = 25931.25
$ws_ARM.Range("L45").Value = 1845.5333
$ws_ARM.Range("M45").Value = -25554.25
$ws_ARM.Range("N45").Value = -2599.5333

# row 61
$ws_ARM.Range("H61").Value = 1384.1923
$ws_ARM.Range("I61").Value = 1285.238
$ws_ARM.Range("K61").Value = 1285.238
$ws_ARM.Range("M61").Value = -1073.238

# row 74
$ws_ARM.Range("H74").Value = 1511.0513
$ws_ARM.Range("I74").Value = 1373.9333
$ws_ARM.Range("J74").Value = 1968.1111
$ws_ARM.Range("K74").Value = 1373.9333
$ws_ARM.Range("L74").Value = 1968.1111
$ws_ARM.Range("M74").Value = -499.9332999999999
$ws_ARM.Range("N74").Value = -3716.1111

# row 77
$ws_ARM.Range("H77").Value = 1511.0513
$ws_ARM.Range("I77").Value = 1373.9333
$ws_ARM.Range("J77").Value = 1968.1111
$ws_ARM.Range("K77").Value = 6869.666499999999
$ws_ARM.Range("L77").Value = 9840.5555
$ws_ARM.Range("M77").Value = -2501.666499999999
$ws_ARM.Range("N77").Value = -18576.5555

# row 132
$ws_ARM.Range("H132").Value = 4482.863
$ws_ARM.Range("I132").Value = 1360.1842
$ws_ARM.Range("J132").Value = 13610.692
$ws_ARM.Range("K132").Value = 4080.5526
$ws_ARM.Range("L132").Value = 40832.076
$ws_ARM.Range("M132").Value = -1550.5526
$ws_ARM.Range("N132").Value = -45892.076

# row 136
$ws_ARM.Range("H136").Value = 1384.1923
$ws_ARM.Range("I136").Value = 1285.238
$ws_ARM.Range("K136").Value = 3855.714
$ws_ARM.Range("M136").Value = -1305.714

# --- BSM ---
# row 80
$ws_BSM.Range("H80").Value = 297.3
$ws_BSM.Range("I80").Value = 56.6
$ws_BSM.Range("J80").Value = 377.53333
$ws_BSM.Range("K80").Value = 56.6
$ws_BSM.Range("L80").Value = 377.53333
$ws_BSM.Range("M80").Value = 941.4
$ws_BSM.Range("N80").Value = -2373.53333

# row 83
$ws_BSM.Range("H83").Value = 297.3
$ws_BSM.Range("I83").Value = 56.6
$ws_BSM.Range("J83").Value = 377.53333
$ws_BSM.Range("K83").Value = 283
$ws_BSM.Range("L83").Value = 1887.66665
$ws_BSM.Range("M83").Value = 4709
$ws_BSM.Range("N83").Value = -11871.66665

# row 122
$ws_BSM.Range("H122").Value = 0
$ws_BSM.Range("J122").Value = 0
$ws_BSM.Range("L122").Value = 0
$ws_BSM.Range("N122").ClearContents()

# row 126
$ws_BSM.Range("H126").Value = 47000
$ws_BSM.Range("J126").Value = 47000
$ws_BSM.Range("L126").Value = 47000
$ws_BSM.Range("N126").Value = -56880

# row 127
$ws_BSM.Range("H127").Value = 80000
$ws_BSM.Range("J127").Value = 80000
$ws_BSM.Range("L127").Value = 80000
$ws_BSM.Range("N127").Value = -89920

# row 130
$ws_BSM.Range("H130").Value = 0
$ws_BSM.Range("J130").Value = 0
$ws_BSM.Range("L130").Value = 0
$ws_BSM.Range("N130").ClearContents()

# row 131
$ws_BSM.Range("H131").Value = 0
$ws_BSM.Range("J131").Value = 0
$ws_BSM.Range("L131").Value = 0
$ws_BSM.Range("N131").ClearContents()

# --- CRP ---
# row 132
$ws_CRP.Range("H132").Value = 2494.0952
$ws_CRP.Range("I132").Value = 2242.8667
$ws_CRP.Range("J132").Value = 3122.1667
$ws_CRP.Range("K132").Value = 6728.6001
$ws_CRP.Range("L132").Value = 9366.500100000001
$ws_CRP.Range("M132").Value = -4198.6001
$ws_CRP.Range("N132").Value = -14426.5001

# --- CUL ---
# row 92
$ws_CUL.Range("H92").Value = 503.8
$ws_CUL.Range("J92").Value = 503.8
$ws_CUL.Range("L92").Value = 1511.4
$ws_CUL.Range("N92").Value = -4007.4

# row 123
$ws_CUL.Range("H123").Value = 8090.4287
$ws_CUL.Range("I123").Value = 3000
$ws_CUL.Range("J123").Value = 8938.833000000001
$ws_CUL.Range("K123").Value = 9000
$ws_CUL.Range("L123").Value = 26816.499
$ws_CUL.Range("M123").Value = -6550
$ws_CUL.Range("N123").Value = -31716.499

# row 129
$ws_CUL.Range("H129").Value = 1705.7059
$ws_CUL.Range("I129").Value = 1383.2222
$ws_CUL.Range("J129").Value = 2068.5
$ws_CUL.Range("K129").Value = 4149.6666
$ws_CUL.Range("L129").Value = 6205.5
$ws_CUL.Range("M129").Value = 850.3334000000004
$ws_CUL.Range("N129").Value = -16205.5

# row 131
$ws_CUL.Range("H131").Value = 24445292
$ws_CUL.Range("I131").Value = 8333805
$ws_CUL.Range("J131").Value = 30304014
$ws_CUL.Range("K131").Value = 25001415
$ws_CUL.Range("L131").Value = 90912042
$ws_CUL.Range("M131").Value = -24996375
$ws_CUL.Range("N131").Value = -90922122

# row 132
$ws_CUL.Range("H132").Value = 1823461.1
$ws_CUL.Range("I132").Value = 2321.6
$ws_CUL.Range("J132").Value = 1986062.9
$ws_CUL.Range("K132").Value = 20894.4
$ws_CUL.Range("L132").Value = 17874566.1
$ws_CUL.Range("M132").Value = -18364.4
$ws_CUL.Range("N132").Value = -17879626.1

# row 137
$ws_CUL.Range("H137").Value = 33345106
$ws_CUL.Range("I137").Value = 12082
$ws_CUL.Range("J137").Value = 50011616
$ws_CUL.Range("K137").Value = 36246
$ws_CUL.Range("L137").Value = 150034848
$ws_CUL.Range("M137").Value = -31146
$ws_CUL.Range("N137").Value = -150045048

# --- GSM ---
# row 122
$ws_GSM.Range("H122").Value = 3829003.5
$ws_GSM.Range("I122").Value = 4052538
$ws_GSM.Range("J122").Value = 3573535.5
$ws_GSM.Range("K122").Value = 12157614
$ws_GSM.Range("L122").Value = 10720606.5
$ws_GSM.Range("M122").Value = -12155164
$ws_GSM.Range("N122").Value = -10725506.5
